$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete TRB column (J, col index 10) - all columns to the right shift left
$ws.Columns.Item(10).Delete()
# Delete PTS column (originally P, now at col index 15 after the first delete)
$ws.Columns.Item(15).Delete()

# Rename last remaining column header (was Playoff_Birth) to Conf
$ws.Cells.Item(1, 16).Value = "Conf"

# Conference (West/East) for each of the 30 teams, rows 2-31
$conf = @("West","West","West","East","East","West","East","West","West","East","West","East","West","West","East","East","West","West","East","East","East","East","East","West","East","East","East","West","West","West")

for ($i = 0; $i -lt 30; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 16).Value = $conf[$i]
}

# Update selection to match the target state
$ws.Range("P32").Select()
